$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("İş Takip Listesi")

$ws1.Range("J2").Formula = "'2025-06-21"
$ws1.Range("K2").Formula = "'2025-11-21"

$ws1.Range("J3").Formula = "'2025-06-21"
$ws1.Range("K3").Formula = "'2025-11-21"

$ws1.Range("J4").Formula = "'2025-06-21"
$ws1.Range("K4").Formula = "'2025-11-21"

$ws1.Range("J5").Formula = "'2025-06-21"
$ws1.Range("K5").Formula = "'2025-11-21"

$ws1.Range("J6").Formula = "'2025-06-21"
$ws1.Range("K6").Formula = "'2025-11-21"

$ws1.Range("J7").Formula = "'2025-06-21"
$ws1.Range("K7").Formula = "'2025-11-21"

$ws1.Range("J8").Formula = "'2025-06-21"
$ws1.Range("K8").Formula = "'2025-11-21"

$ws1.Range("J9").Formula = "'2025-06-21"
$ws1.Range("K9").Formula = "'2025-11-21"

$ws1.Range("J10").Formula = "'2025-06-21"
$ws1.Range("K10").Formula = "'2025-11-21"

$ws1.Range("J33").Formula = "'2025-06-23"
$ws1.Range("K33").Formula = "'2025-11-23"

$ws1.Range("J34").Formula = "'2025-06-23"
$ws1.Range("K34").Formula = "'2025-11-23"

$ws1.Range("J35").Formula = "'2025-06-23"
$ws1.Range("K35").Formula = "'2025-11-23"

$ws1.Range("J36").Formula = "'2025-06-23"
$ws1.Range("K36").Formula = "'2025-11-23"

$ws1.Range("J37").Formula = "'2025-06-23"
$ws1.Range("K37").Formula = "'2025-11-23"

$ws1.Range("J38").Formula = "'2025-06-23"
$ws1.Range("K38").Formula = "'2025-11-23"

$ws1.Range("J39").Formula = "'2025-06-23"
$ws1.Range("K39").Formula = "'2025-11-23"

$ws1.Range("J40").Formula = "'2025-06-23"
$ws1.Range("K40").Formula = "'2025-11-23"

$ws1.Range("J41").Formula = "'2025-06-23"
$ws1.Range("K41").Formula = "'2025-11-23"

$ws1.Range("J42").Formula = "'2025-06-23"
$ws1.Range("K42").Formula = "'2025-11-23"

$ws1.Range("J43").Formula = "'2025-06-23"
$ws1.Range("K43").Formula = "'2025-11-23"

$ws1.Range("J44").Formula = "'2025-06-23"
$ws1.Range("K44").Formula = "'2025-11-23"

$ws1.Range("J45").Formula = "'2025-06-23"
$ws1.Range("K45").Formula = "'2025-11-23"

$ws1.Range("J46").Formula = "'2025-06-23"
$ws1.Range("K46").Formula = "'2025-11-23"

$ws1.Range("J47").Formula = "'2025-06-23"
$ws1.Range("K47").Formula = "'2025-11-23"

$ws1.Range("J48").Formula = "'2025-06-23"
$ws1.Range("K48").Formula = "'2025-11-23"

$ws1.Range("J49").Formula = "'2025-06-23"
$ws1.Range("K49").Formula = "'2025-11-23"

$ws1.Range("J50").Formula = "'2025-06-23"
$ws1.Range("K50").Formula = "'2025-11-23"

$ws1.Range("J51").Formula = "'2025-06-23"
$ws1.Range("K51").Formula = "'2025-11-23"

$ws1.Range("J52").Formula = "'2025-06-23"
$ws1.Range("K52").Formula = "'2025-11-23"

$ws1.Range("J53").Formula = "'2025-06-23"
$ws1.Range("K53").Formula = "'2025-11-23"

$ws1.Range("J54").Formula = "'2025-06-23"
$ws1.Range("K54").Formula = "'2025-11-23"

$ws1.Range("J55").Formula = "'2025-06-23"
$ws1.Range("K55").Formula = "'2025-11-23"

$ws1.Range("J56").Formula = "'2025-06-23"
$ws1.Range("K56").Formula = "'2025-11-23"

$ws1.Range("J57").Formula = "'2025-06-23"
$ws1.Range("K57").Formula = "'2025-11-23"

$ws1.Range("J58").Formula = "'2025-06-23"
$ws1.Range("K58").Formula = "'2025-11-23"

$ws1.Range("J59").Formula = "'2025-06-23"
$ws1.Range("K59").Formula = "'2025-11-23"

$ws1.Range("J60").Formula = "'2025-06-23"
$ws1.Range("K60").Formula = "'2025-11-23"

$ws1.Range("J61").Formula = "'2025-06-23"
$ws1.Range("K61").Formula = "'2025-11-23"

$ws1.Range("J62").Formula = "'2025-06-23"
$ws1.Range("K62").Formula = "'2025-11-23"

$ws1.Range("J63").Formula = "'2025-06-23"
$ws1.Range("K63").Formula = "'2025-11-23"

$ws1.Range("J64").Formula = "'2025-06-23"
$ws1.Range("K64").Formula = "'2025-11-23"

$ws1.Range("J65").Formula = "'2025-06-23"
$ws1.Range("K65").Formula = "'2025-11-23"

$ws1.Range("J66").Formula = "'2025-06-23"
$ws1.Range("K66").Formula = "'2025-11-23"

$ws1.Range("J67").Formula = "'2025-06-23"
$ws1.Range("K67").Formula = "'2025-11-23"

$ws1.Range("J68").Formula = "'2025-06-23"
$ws1.Range("K68").Formula = "'2025-11-23"

$ws1.Range("J69").Formula = "'2025-06-23"
$ws1.Range("K69").Formula = "'2025-11-23"

$ws1.Range("J70").Formula = "'2025-06-23"
$ws1.Range("K70").Formula = "'2025-11-23"

$ws1.Range("J71").Formula = "'2025-06-23"
$ws1.Range("K71").Formula = "'2025-11-23"

$ws1.Range("J72").Formula = "'2025-06-23"
$ws1.Range("K72").Formula = "'2025-11-23"
$ws1.Range("L72").Value = 'ASKI İLANINA HAZIRLANIYOR'

$ws1.Range("J73").Formula = "'2025-06-23"
$ws1.Range("K73").Formula = "'2025-11-23"

$ws1.Range("J74").Formula = "'2025-06-23"
$ws1.Range("K74").Formula = "'2025-11-23"

$ws1.Range("J75").Formula = "'2025-06-23"
$ws1.Range("K75").Formula = "'2025-11-23"

$ws1.Range("J76").Formula = "'2025-06-23"
$ws1.Range("K76").Formula = "'2025-11-23"

$ws1.Range("J77").Formula = "'2025-06-23"
$ws1.Range("K77").Formula = "'2025-11-23"

$ws1.Range("J78").Formula = "'2025-06-23"
$ws1.Range("K78").Formula = "'2025-11-23"

$ws1.Range("J79").Formula = "'2025-06-23"
$ws1.Range("K79").Formula = "'2025-11-23"

$ws1.Range("J80").Formula = "'2025-06-23"
$ws1.Range("K80").Formula = "'2025-11-23"

$ws1.Range("J81").Formula = "'2025-06-23"
$ws1.Range("K81").Formula = "'2025-11-23"

$ws1.Range("J82").Formula = "'2025-06-23"
$ws1.Range("K82").Formula = "'2025-11-23"

$ws1.Range("J83").Formula = "'2025-06-23"
$ws1.Range("K83").Formula = "'2025-11-23"

$ws1.Range("J84").Formula = "'2025-06-23"
$ws1.Range("K84").Formula = "'2025-11-23"

$ws1.Range("J85").Formula = "'2025-06-23"
$ws1.Range("K85").Formula = "'2025-11-23"

$ws1.Range("J86").Formula = "'2025-06-23"
$ws1.Range("K86").Formula = "'2025-11-23"

$ws1.Range("J87").Formula = "'2025-06-23"
$ws1.Range("K87").Formula = "'2025-11-23"

$ws1.Range("J88").Formula = "'2025-06-23"
$ws1.Range("K88").Formula = "'2025-11-23"

$ws1.Range("J89").Formula = "'2025-06-23"
$ws1.Range("K89").Formula = "'2025-11-23"

$ws1.Range("J90").Formula = "'2025-06-23"
$ws1.Range("K90").Formula = "'2025-11-23"

$ws1.Range("J91").Formula = "'2025-06-23"
$ws1.Range("K91").Formula = "'2025-11-23"

$ws1.Range("J92").Formula = "'2025-06-23"
$ws1.Range("K92").Formula = "'2025-11-23"

$ws1.Range("J93").Formula = "'2025-06-23"
$ws1.Range("K93").Formula = "'2025-11-23"

$ws1.Range("J94").Formula = "'2025-06-23"
$ws1.Range("K94").Formula = "'2025-11-23"

$ws1.Range("J95").Formula = "'2024-04-21"
$ws1.Range("K95").Formula = "'2025-06-15"

$ws1.Range("J96").Formula = "'2024-04-21"
$ws1.Range("K96").Formula = "'2025-06-15"

$ws1.Range("J97").Formula = "'2024-04-21"
$ws1.Range("K97").Formula = "'2025-06-15"

$ws1.Range("J98").Formula = "'2024-04-21"
$ws1.Range("K98").Formula = "'2025-06-15"

$ws1.Range("J99").Formula = "'2024-04-21"
$ws1.Range("K99").Formula = "'2025-06-15"

$ws1.Range("J100").Formula = "'2024-04-21"
$ws1.Range("K100").Formula = "'2025-06-15"

$ws1.Range("J101").Formula = "'2024-04-21"
$ws1.Range("K101").Formula = "'2025-06-15"

$ws1.Range("J102").Formula = "'2024-04-21"
$ws1.Range("K102").Formula = "'2025-06-15"

$ws1.Range("J103").Formula = "'2024-04-21"
$ws1.Range("K103").Formula = "'2025-06-15"

$ws1.Range("J104").Formula = "'2024-04-21"
$ws1.Range("K104").Formula = "'2025-06-15"

$ws1.Range("J105").Formula = "'2024-04-21"
$ws1.Range("K105").Formula = "'2025-06-15"

$ws1.Range("J106").Formula = "'2024-04-21"
$ws1.Range("K106").Formula = "'2025-06-15"

$ws1.Range("J107").Formula = "'2024-04-21"
$ws1.Range("K107").Formula = "'2025-06-15"

$ws1.Range("J108").Formula = "'2024-04-21"
$ws1.Range("K108").Formula = "'2025-06-15"

$ws1.Range("J109").Formula = "'2024-04-21"
$ws1.Range("K109").Formula = "'2025-06-15"

$ws1.Range("J110").Formula = "'2024-04-21"
$ws1.Range("K110").Formula = "'2025-06-15"

$ws1.Range("J111").Formula = "'2024-04-21"
$ws1.Range("K111").Formula = "'2025-06-15"

$ws1.Range("J112").Formula = "'2024-04-21"
$ws1.Range("K112").Formula = "'2025-06-15"

$ws1.Range("J113").Formula = "'2024-04-21"
$ws1.Range("K113").Formula = "'2025-06-15"

$ws1.Range("J114").Formula = "'2024-04-21"
$ws1.Range("K114").Formula = "'2025-06-15"

$ws1.Range("J115").Formula = "'2024-04-21"
$ws1.Range("K115").Formula = "'2025-06-15"

$ws1.Range("J116").Formula = "'2024-04-21"
$ws1.Range("K116").Formula = "'2025-06-15"

$ws1.Range("J117").Formula = "'2024-04-21"
$ws1.Range("K117").Formula = "'2025-06-15"

$ws1.Range("J118").Formula = "'2024-04-21"
$ws1.Range("K118").Formula = "'2025-06-15"

$ws1.Range("J119").Formula = "'2024-04-21"
$ws1.Range("K119").Formula = "'2025-06-15"

$ws1.Range("J120").Formula = "'2024-04-21"
$ws1.Range("K120").Formula = "'2025-06-15"

$ws1.Range("J121").Formula = "'2024-04-21"
$ws1.Range("K121").Formula = "'2025-06-15"

$ws1.Range("J122").Formula = "'2024-04-21"
$ws1.Range("K122").Formula = "'2025-06-15"

$ws2 = $wb.Worksheets.Item("Güncelleme")

$ws2.Range("J2").Formula = "'2024-07-28"
$ws2.Range("N2").Formula = "'2025-03-31"
$ws2.Range("P2").Formula = "'2025-06-17"

$ws2.Range("J3").Formula = "'2024-10-29"
$ws2.Range("N3").Formula = "'2025-07-19"
$ws2.Range("P3").Formula = "'2025-10-27"

$ws2.Range("J4").Formula = "'2024-09-02"
$ws2.Range("N4").Formula = "'2025-02-23"
$ws2.Range("P4").Formula = "'2025-05-20"

$ws2.Range("I5").Formula = "'2025-02-26"

$ws2.Range("J6").Formula = "'2025-10-09"
$ws2.Range("N6").Formula = "'2025-06-28"
$ws2.Range("P6").Formula = "'2025-12-24"

$ws2.Range("I7").Formula = "'2024-10-29"
$ws2.Range("J7").Formula = "'2024-10-29"

$ws2.Range("J8").Formula = "'2024-10-17"
$ws2.Range("N8").Formula = "'2025-03-18"
$ws2.Range("P8").Formula = "'2025-04-20"

$ws2.Range("I9").Formula = "'2025-06-14"
$ws2.Range("J9").Formula = "'2024-11-30"

$ws2.Range("J10").Formula = "'2024-09-28"
$ws2.Range("N10").Formula = "'2025-07-09"
$ws2.Range("P10").Formula = "'2025-10-31"

$ws2.Range("I11").Formula = "'2025-04-05"
$ws2.Range("J11").Formula = "'2024-11-11"
$ws2.Range("N11").Formula = "'2025-07-29"
$ws2.Range("P11").Formula = "'2025-12-24"

$ws2.Range("J12").Formula = "'2024-10-09"
$ws2.Range("N12").Formula = "'2025-06-18"
$ws2.Range("P12").Formula = "'2025-10-21"

$ws2.Range("J13").Formula = "'2024-12-07"

$ws2.Range("J14").Formula = "'2025-10-05"
$ws2.Range("N14").Formula = "'2025-10-27"

$ws2.Range("J15").Formula = "'2024-12-26"
$ws2.Range("N15").Formula = "'2025-07-16"
$ws2.Range("P15").Formula = "'2025-10-28"

$ws2.Range("J16").Formula = "'2024-08-24"
$ws2.Range("N16").Formula = "'2025-02-01"
$ws2.Range("P16").Formula = "'2025-04-20"

$ws2.Range("J17").Formula = "'2024-09-09"
$ws2.Range("N17").Formula = "'2025-10-27"

$ws2.Range("J18").Formula = "'2025-02-16"

$ws2.Range("I19").Formula = "'2025-04-06"
$ws2.Range("J19").Formula = "'2024-12-26"
$ws2.Range("N19").Formula = "'2025-08-05"

$ws2.Range("J20").Formula = "'2024-12-07"
$ws2.Range("N20").Formula = "'2025-12-11"

$ws2.Range("J21").Formula = "'2024-09-30"

$ws2.Range("J22").Formula = "'2024-09-30"

$ws2.Range("J23").Formula = "'2024-12-08"

$ws2.Range("I24").Formula = "'2025-06-04"

$ws2.Range("J25").Formula = "'2024-11-02"

$ws2.Range("J27").Formula = "'2025-01-23"

$ws2.Range("J28").Formula = "'2024-11-21"
$ws2.Range("N28").Formula = "'2025-11-06"

$ws2.Range("I29").Formula = "'2025-02-10"
$ws2.Range("J29").Formula = "'2024-12-08"
$ws2.Range("N29").Formula = "'2025-10-18"
